$wb = $excel.ActiveWorkbook

# --- Create the two new sheets, in tab order ColumnHeadersEims, ColumnHeadersToi, before CategoricalVariables ---
$catRef = $wb.Worksheets.Item("CategoricalVariables")
$ws1 = $wb.Worksheets.Add($catRef)
$ws1.Name = "ColumnHeadersEims"

$catRef2 = $wb.Worksheets.Item("CategoricalVariables")
$ws2 = $wb.Worksheets.Add($catRef2)
$ws2.Name = "ColumnHeadersToi"

# --- Populate ColumnHeadersEims (sheet1) ---
$ws1.Range('A1').Value = 'attributeName'
$ws1.Range('B1').Value = 'attributeDefinition'
$ws1.Range('C1').Value = 'class'
$ws1.Range('D1').Value = 'unit'
$ws1.Range('E1').Value = 'dateTimeFormatString'
$ws1.Range('F1').Value = 'missingValueCode'
$ws1.Range('G1').Value = 'missingValueCodeExplanation'
$ws1.Range('A2').Value = 'cruise'
$ws1.Range('B2').Value = 'Identifier for research cruise generally including abbreviation for research vessel and voyage number'
$ws1.Range('C2').Value = 'character'
$ws1.Range('A3').Value = 'datetime_utc_matlab'
$ws1.Range('B3').Value = 'PI-provided UTC date and time '
$ws1.Range('C3').Value = 'Date'
$ws1.Range('E3').Value = 'YYYY-MM-DD hh:mm:ss'
$ws1.Range('A4').Value = 'latitude_matlab'
$ws1.Range('B4').Value = 'Latitude of sample event provided by PI'
$ws1.Range('C4').Value = 'numeric'
$ws1.Range('D4').Value = 'degree'
$ws1.Range('A5').Value = 'longitude_matlab'
$ws1.Range('B5').Value = 'Longitude of sample event provided by PI'
$ws1.Range('C5').Value = 'numeric'
$ws1.Range('D5').Value = 'degree'
$ws1.Range('A6').Value = 'depth'
$ws1.Range('B6').Value = 'Data product depth of sample below sea surface, for underway samples depth of ship''s intake. URI http://vocab.nerc.ac.uk/collection/P09/current/DEPH/'
$ws1.Range('C6').Value = 'numeric'
$ws1.Range('D6').Value = 'meter'
$ws1.Range('A7').Value = 'biosat'
$ws1.Range('B7').Value = 'Percent biological saturation, [ (O2/Ar) meas / (O2/Ar)equilibrium - 1 ]x100'
$ws1.Range('C7').Value = 'numeric'
$ws1.Range('D7').Value = 'dimensionless'
$ws1.Range('F7').Value = 'NaN'
$ws1.Range('G7').Value = 'Missing value'
$ws1.Range('A8').Value = 'O2_Ar_ratio'
$ws1.Range('B8').Value = 'Oxygen-argon ratio of EIMS sample from underway corrected for air values'
$ws1.Range('C8').Value = 'numeric'
$ws1.Range('D8').Value = 'dimensionless'
$ws1.Range('F8').Value = 'NaN'
$ws1.Range('G8').Value = 'Missing value'

# --- Populate ColumnHeadersToi (sheet2) ---
$ws2.Range('A1').Value = 'attributeName'
$ws2.Range('B1').Value = 'attributeDefinition'
$ws2.Range('C1').Value = 'class'
$ws2.Range('D1').Value = 'unit'
$ws2.Range('E1').Value = 'dateTimeFormatString'
$ws2.Range('F1').Value = 'missingValueCode'
$ws2.Range('G1').Value = 'missingValueCodeExplanation'
$ws2.Range('A2').Value = 'cruise'
$ws2.Range('B2').Value = 'Identifier for research cruise generally including abbreviation for research vessel and voyage number'
$ws2.Range('C2').Value = 'character'
$ws2.Range('A3').Value = 'datetime_utc'
$ws2.Range('B3').Value = 'Data product UTC date and time'
$ws2.Range('C3').Value = 'Date'
$ws2.Range('E3').Value = 'YYYY-MM-DD hh:mm:ss'
$ws2.Range('A4').Value = 'datetime_utc_matlab'
$ws2.Range('B4').Value = 'PI-provided UTC date and time'
$ws2.Range('C4').Value = 'Date'
$ws2.Range('E4').Value = 'YYYY-MM-DD hh:mm:ss'
$ws2.Range('A5').Value = 'latitude_API'
$ws2.Range('B5').Value = 'Latitude of sample event provided by NES-LTER API'
$ws2.Range('C5').Value = 'numeric'
$ws2.Range('D5').Value = 'degree'
$ws2.Range('A6').Value = 'longitude_API'
$ws2.Range('B6').Value = 'Longitude of sample event provided by NES-LTER API'
$ws2.Range('C6').Value = 'numeric'
$ws2.Range('D6').Value = 'degree'
$ws2.Range('A7').Value = 'toi_source'
$ws2.Range('B7').Value = 'Source of bottle sample whether from Niskin or underway'
$ws2.Range('C7').Value = 'categorical'
$ws2.Range('A8').Value = 'cast'
$ws2.Range('B8').Value = 'CTD rosette cast number, chronological per cruise'
$ws2.Range('C8').Value = 'numeric'
$ws2.Range('D8').Value = 'dimensionless'
$ws2.Range('F8').Value = 'NaN'
$ws2.Range('G8').Value = 'Sample from underway'
$ws2.Range('A9').Value = 'niskin'
$ws2.Range('B9').Value = 'Rosette bottle position number'
$ws2.Range('C9').Value = 'numeric'
$ws2.Range('D9').Value = 'dimensionless'
$ws2.Range('F9').Value = 'NaN'
$ws2.Range('G9').Value = 'Sample from underway'
$ws2.Range('A10').Value = 'depth'
$ws2.Range('B10').Value = 'Data product depth of sample below sea surface, for underway samples depth of ship''s intake, for Niskins from CTD summary data in NES-LTER  API. URI http://vocab.nerc.ac.uk/collection/P09/current/DEPH/'
$ws2.Range('C10').Value = 'numeric'
$ws2.Range('D10').Value = 'meter'
$ws2.Range('A11').Value = 'depth_matlab'
$ws2.Range('B11').Value = 'PI-provided depth of sample below sea surface. '
$ws2.Range('C11').Value = 'numeric'
$ws2.Range('D11').Value = 'meter'
$ws2.Range('A12').Value = 'O2_Ar_delta'
$ws2.Range('B12').Value = 'Oxygen-argon ratio divided by the reference ratio (oxygen-argon ratio in air minus 1, multiplied by 100)'
$ws2.Range('C12').Value = 'numeric'
$ws2.Range('D12').Value = 'dimensionless'
$ws2.Range('A13').Value = 'O2_Ar_ratio'
$ws2.Range('B13').Value = 'Oxygen-argon ratio of bottle sample'
$ws2.Range('C13').Value = 'numeric'
$ws2.Range('D13').Value = 'dimensionless'
$ws2.Range('A14').Value = 'cap_Delta_17O'
$ws2.Range('B14').Value = 'Triple isotopic composition of dissolved oxygen versus atmospheric O2, D17O'
$ws2.Range('C14').Value = 'numeric'
$ws2.Range('D14').Value = 'perMeg'
$ws2.Range('A15').Value = 'd17O'
$ws2.Range('B15').Value = 'Enrichment of oxygen-17 in dissolved oxygen (delta(17)O) in the water body by mass spectrometry'
$ws2.Range('C15').Value = 'numeric'
$ws2.Range('D15').Value = 'perMil'
$ws2.Range('A16').Value = 'd18O'
$ws2.Range('B16').Value = 'Enrichment of oxygen-18 in dissolved oxygen {18O in O2 CAS 14797-71-8} {delta(18)O} in the water body by mass spectrometry URI http://vocab.nerc.ac.uk/collection/P01/current/D18OMXDG/'
$ws2.Range('C16').Value = 'numeric'
$ws2.Range('D16').Value = 'perMil'

# --- Formatting: wrap-text header/definition cells ---
$ws1.Range("B1").WrapText = $true
$ws1.Range("B6").WrapText = $true
$ws1.Rows.Item(6).RowHeight = 31.2

$ws2.Range("B1").WrapText = $true
$ws2.Range("B10").WrapText = $true
$ws2.Rows.Item(10).RowHeight = 46.8

# --- Special pasted-in font color for the d17O / d18O definitions ---
$ws2.Range("B15").Font.Color = 1907741
$ws2.Range("B16").Font.Color = 1907741

# --- Column widths ---
$ws1.Columns.Item(1).ColumnWidth = 24.33
$ws1.Columns.Item(2).ColumnWidth = 76.67

$ws2.Columns.Item(1).ColumnWidth = 15.67
$ws2.Columns.Item(2).ColumnWidth = 80.67

# --- Print setup for ColumnHeadersToi ---
$ws2.PageSetup.Orientation = 1

# --- Selections matching the authored file ---
$ws1.Range("B29").Select() | Out-Null
$ws2.Range("B29").Select() | Out-Null

# --- Leave ColumnHeadersToi as the active sheet/tab (matches activeTab=1) ---
$ws2.Activate()
$ws2.Range("B29").Select() | Out-Null
